$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the acquisition timestamp in column A for rows 2-15 to reflect the
# latest append run at 2025-11-11 12:48:49.
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-11-11 12:48:49"
}
